$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (count) values for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (updated meanEMG / legmaxROM values) for columns B:E
$ws.Range("B2").Value = 12.868937944505966
$ws.Range("C2").Value = 11.851118303041785
$ws.Range("D2").Value = 13.348431705634853
$ws.Range("E2").Value = 12.782451313437109

# Row 3 (updated meanEMG / legmaxROM values) for columns B:E
$ws.Range("B3").Value = 12.930522173316671
$ws.Range("C3").Value = 10.883366192557336
$ws.Range("D3").Value = 14.630921403483496
$ws.Range("E3").Value = 11.881610585812835

# Match the narrowed selection recorded in the saved workbook
$ws.Range("B1:E3").Select()
